# "Editing my personal information" — update the fund-number (Numero) column
# to numeric type for existing rows, and append the newest reporting period
# (2024-04-30, serial 45412) for the four funds.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix C49:C52: the fund account numbers were stored as text; store them
#     as real numbers instead (same visible digits). ---
$ws.Cells.Item(49, 3).Value = 252000001274
$ws.Cells.Item(50, 3).Value = 1111000544148
$ws.Cells.Item(51, 3).Value = 252000011589
$ws.Cells.Item(52, 3).Value = 342000006519

# --- Append new rows 53-56 for period 2024-04-30 (serial 45412) ---
$rows = @(
    @{ r=53; b="Fidurenta";        c="1111000544148"; d=45233.36324127; e=119590.77;         f=0;       g=0;       h=-133.75;    i=0;                  j=119457.02;  k=1.35  },
    @{ r=54; b="Renta Acciones";   c="252000001274";  d=70711.08490446; e=779626.88;         f=6000000; g=0;       h=-147388.81; i=0;                  j=6632238.07; k=20.62 },
    @{ r=55; b="Renta Fija Plazo"; c="252000011589";  d=35671.28200578; e=9256619.279999999; f=0;       g=6000000; h=-130131.73; i=2871.7;             j=3123615.85; k=13.57 },
    @{ r=56; b="Fiducuenta";       c="342000006519";  d=39507.92241913; e=596031.51;         f=2404134; g=1589059; h=8401.42;    i=74.31999999999999; j=1419433.61; k=6.7   }
)

# Reference cell for "plain/default" formatting (never written -> style 0),
# used to strip the quote-prefix format that typing a leading apostrophe adds.
$blankFormatCell = $ws.Cells.Item(200, 26)

foreach ($row in $rows) {
    $r = $row.r

    # Column A: date, formatted like the rows above it (style copied from A49).
    $ws.Cells.Item($r, 1).Value = 45412
    $ws.Cells.Item(49, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 2).Value = $row.b

    # Column C: fund number kept as TEXT for the new rows (leading apostrophe
    # forces text-type entry), then strip the resulting quote-prefix style.
    $ws.Cells.Item($r, 3).Value = "'" + $row.c
    $blankFormatCell.Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)

    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
    $ws.Cells.Item($r, 8).Value = $row.h
    $ws.Cells.Item($r, 9).Value = $row.i
    $ws.Cells.Item($r, 10).Value = $row.j
    $ws.Cells.Item($r, 11).Value = $row.k
}
